$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings used on row 1 (and wherever else) so that
# "<x>_old" -> "<x>_FV2210" and "<x>_new" -> "<x>_FV2304"
$words = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            foreach ($w in $words) {
                if ($val -eq ($w + "_old")) {
                    $cell.Value2 = $w + "_FV2210"
                    $val = $cell.Value2
                } elseif ($val -eq ($w + "_new")) {
                    $cell.Value2 = $w + "_FV2304"
                    $val = $cell.Value2
                }
            }
        }
    }
}

# Freeze the header row (row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a table (ListObject)
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
